$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B6").Value = "HE130577"
$ws.Range("I16").Select()
